$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "From Query" example table (D7:D11) originally stored the same
# date-time stamp (40720.2728356481) in every row. Replace it with five
# distinct, date-only values (2000-01-01 .. 2000-01-05), written as raw
# serial numbers so the existing date style isn't disturbed by date-string
# auto-formatting.
$ws.Range("D7").Value = 36526
$ws.Range("D8").Value = 36527
$ws.Range("D9").Value = 36528
$ws.Range("D10").Value = 36529
$ws.Range("D11").Value = 36530

# Switch the column's display format from date+time ("m/d/yyyy h:mm",
# numFmtId 22) to date-only ("m/d/yyyy", numFmtId 14).
$ws.Range("D7:D11").NumberFormat = "mm-dd-yy"

# The column is narrower now that it only needs to fit a short date.
$ws.Range("D1").ColumnWidth = 9.5
